$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# --- Header row (row 1): insert "category" before the existing "date" column,
# and append "legislator_id", "source_file", "index" as new trailing columns.
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Data rows: shift H/I/J content to make room for the new "category" column,
# move the legislator name into K (leaving the numeric id in a new L column),
# and append the new source_file / index columns.
$rows = @(2, 3, 4)
$indexValues = @(90, 91, 92)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]

    $ws.Cells.Item($r, 8).Value = "stock"          # H: property_category
    $ws.Cells.Item($r, 9).Value = "normal"         # I: category (new)
    $ws.Cells.Item($r, 10).Value = "2012-05-01"    # J: date
    $ws.Cells.Item($r, 11).Value = "翁重鈞"         # K: legislator_name
    $ws.Cells.Item($r, 12).Value = 551             # L: legislator_id (new position)
    $ws.Cells.Item($r, 13).Value = "tmp6aad1"      # M: source_file (new)
    $ws.Cells.Item($r, 14).Value = $indexValues[$i] # N: index (new)
}
